# Scene.xlsx fix-up:
#   - "FilePath" column (F9:F14) pointed at "../resource/Ini/Scene/N.xml" with
#     a mixed-case "Ini" directory; the real resource folder on disk is
#     lower-case "ini", so every reference is corrected in place.
#   - The sheet's saved cursor position moves from F16 to G17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataNode")

# FilePath values live in column F, rows 9-14 (one row per scene record).
# Only the "Ini" -> "ini" path segment changes; the rest of each path is
# left untouched.
$filePathCol = 6
$firstRow = 9
$lastRow = 14

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $filePathCol)
    $current = $cell.Value2
    if ($current -ne $null -and $current -like "*/Ini/Scene/*") {
        $cell.Value = $current -replace "/Ini/Scene/", "/ini/Scene/"
    }
}

# Restore the sheet's active-cell selection to G17 (was F16).
$ws.Range("G17").Select()
